$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,(1,1,"backside flip",$null,1,$null,"flat","curb",$null,"Union - Right to Skate")
    ,(1,2,"switch varial",1,1,$null,"flat","curb",$null,"Union - Right to Skate")
    ,(2,3,"fs 5-0 to switch crook",$null,$null,$null,"ledge","curb",$null,"Union - Right to Skate")
    ,(3,4,"fs 180",$null,1,$null,"flat","curb",$null,"Union - Right to Skate")
    ,(3,5,"switch hardflip",1,1,$null,"flat","curb",$null,"Union - Right to Skate")
    ,(4,6,"switch bs noseslide 270",1,$null,$null,"ledge","curb",$null,"Union - Right to Skate")
    ,(5,7,"halfcab to noseslide to nosegrind to regular",$null,1,$null,"ledge","curb",$null,"Union - Right to Skate")
    ,(5,8,"nollie bs flip",1,1,$null,"flat",$null,$null,"Union - Right to Skate")
    ,(6,9,"bs tailslide",$null,1,$null,"ledge","curb",$null,"Union - Right to Skate")
    ,(6,10,"fs 180 nosegrind to noseslide to fakie",$null,1,$null,"ledge","curb",$null,"Union - Right to Skate")
    ,(6,11,"switch bs 180",1,1,$null,"flat",$null,$null,"Union - Right to Skate")
    ,(7,12,"nollie noseslide to fakie",1,1,$null,"ledge",$null,$null,"Union - Right to Skate")
    ,(7,13,"switch fs tailslide",1,1,$null,"ledge",$null,$null,"Union - Right to Skate")
    ,(8,14,"nollie bs 180 to fakie 5-0",1,$null,$null,"ledge","curb",$null,"Union - Right to Skate")
    ,(9,15,"fs 5-0 to fs noseslide",$null,$null,$null,"ledge",$null,$null,"Union - Right to Skate")
    ,(10,16,"switch bs noseslide 270 popover",1,$null,1,"ledge","curb",$null,"Union - Right to Skate")
    ,(11,17,"switch bs flip",1,$null,$null,"flat","flatgap",$null,"Union - Right to Skate")
    ,(12,18,"switch big flip",1,$null,$null,"stair",$null,$null,"Union - Right to Skate")
    ,(13,19,"bs tailslide to fakie",$null,$null,$null,"ledge","shootout",$null,"Union - Right to Skate")
    ,(14,20,"bs noseslide 270",$null,$null,$null,"ledge","shootout",$null,"Union - Right to Skate")
    ,(15,21,"bs tailslide to fakie 5-0 to fakie",$null,$null,$null,"ledge","curb",$null,"Union - Right to Skate")
    ,(16,22,"nollie fs tailslide",1,$null,$null,"ledge",$null,$null,"Union - Right to Skate")
    ,(17,23,"kickflip bs nosedlide to fakie",$null,$null,$null,"ledge",$null,$null,"Union - Right to Skate")
    ,(18,24,"fs flip to switch manual to regular",$null,$null,$null,"manual",$null,$null,"Union - Right to Skate")
    ,(19,25,"melon grab",$null,$null,$null,"gap",$null,$null,"Union - Right to Skate")
    ,(20,26,"nollie fs 180",1,$null,$null,"gap",$null,"carlsbad","Union - Right to Skate")
    ,(21,27,"bs darkslide",$null,$null,$null,"ledge",$null,$null,"Union - Right to Skate")
    ,(22,28,"bs tailslide to fakie",$null,1,$null,"ledge",$null,$null,"Union - Right to Skate")
    ,(23,29,"halfcab heelflip to noseslide to fakie",$null,1,$null,"ledge",$null,$null,"Union - Right to Skate")
    ,(24,30,"switch flip bs tail",1,$null,$null,"ledge","curb",$null,"Union - Right to Skate")
)

$startRow = 35
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        if ($null -ne $row[$c]) {
            $ws.Cells.Item($r, $c + 1).Value = $row[$c]
        }
    }
}

# Move the active selection to reflect where the user ended up after
# appending the new rows (matches the post-edit cursor position).
$ws.Range("A65").Select()

Write-Host "Added rows 35-64"